# Applies the edit described by the diff:
#  - Clears the URL / username / password test values in row 2 (M2:O2)
#    of the Input_Value sheet (they contained sensitive sample creds that
#    were scrubbed before upload).
#  - Moves the active selection on that sheet to reflect where the user
#    was last working (M2:O2, scrolled so column L is left-most).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate()

# Clear the previously-populated URL / user / password sample values.
$ws.Range("M2:O2").ClearContents()

# Reflect the new selection / scroll position recorded in the workbook.
$ws.Range("L1").Select()
$ws.Range("M2:O2").Select()
$excel.ActiveWindow.ScrollColumn = 12
